$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.035.77'
$ws.Range("E2").Value = '  -0.86%  '

$ws.Range("D3").Value = '2.616.20'
$ws.Range("E3").Value = '  +0.76%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.15'
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.28'
$ws.Range("E6").Value = '  -3.50%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  -4.92%  '

$ws.Range("D9").Value = '2.620.28'
$ws.Range("E9").Value = '  +0.56%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.32'
$ws.Range("E10").Value = '  -5.07%  '

$ws.Range("E11").Value = '  -0.16%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.340'
$ws.Range("E12").Value = '  -1.79%  '

$ws.Range("E13").Value = '  -0.44%  '

$ws.Range("D14").Value = '3.077.61'
$ws.Range("E14").Value = '  +0.83%  '

$ws.Range("D15").Value = '60.073.04'
$ws.Range("E15").Value = '  -0.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.15'
$ws.Range("E16").Value = '  -2.50%  '

$ws.Range("E17").Value = '  -1.87%  '

$ws.Range("D18").Value = '2.615.96'
$ws.Range("E18").Value = '  +0.53%  '

$ws.Range("E19").Value = '  -2.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '339.76'
$ws.Range("E20").Value = '  -3.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.41'
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("E22").Value = '  -1.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.72'
$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.418'
$ws.Range("E25").Value = '  -2.12%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.15%  '

$ws.Range("E27").Value = '  -3.02%  '

$ws.Range("D28").Value = '0.0₃0804'
$ws.Range("E28").Value = '  -4.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.02'
$ws.Range("E29").Value = '  -4.35%  '

$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("E31").Value = '  -0.35%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.95'
$ws.Range("E32").Value = '  -2.20%  '

$ws.Range("B33").Value = 'Aptos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.96'
$ws.Range("E33").Value = '  -5.07%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.91'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.93'
$ws.Range("E35").Value = '  -5.52%  '

$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("E37").Value = '  -5.36%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.858'
$ws.Range("E38").Value = '  +2.25%  '

$ws.Range("E39").Value = '  +0.23%  '

$ws.Range("E40").Value = '  -4.41%  '

$ws.Range("E41").Value = '  -4.41%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '286.15'
$ws.Range("E42").Value = '  +0.09%  '

$ws.Range("E43").Value = '  +0.21%  '

$ws.Range("E44").Value = '  -1.29%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.998'
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0544'
$ws.Range("E46").Value = '  -2.48%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.38'
$ws.Range("E47").Value = '  -0.75%  '

$ws.Range("E48").Value = '  +0.87%  '

$ws.Range("E49").Value = '  -2.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.66'
$ws.Range("E50").Value = '  -3.33%  '

$ws.Range("D51").Value = '1.960.44'
$ws.Range("E51").Value = '  -0.09%  '
